{"js": "// Update the two-digit/one-digit division exercises in the first table.\n// Each data row of the table (rows 0, 4, 8, 12, 16 \u2014 the other rows are\n// blank \"answer\" rows) holds five \"NN\u00f7N=\" expressions; replace them in\n// reading order with their new values, matched positionally (not by text,\n// since some old values repeat with different replacements, e.g. \"43\u00f72=\"\n// appears twice and must become two different results).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// [oldText, newText] pairs, in document (reading) order.\nconst pairs = [\n  [\"30\u00f74=\", \"45\u00f75=\"],\n  [\"11\u00f77=\", \"56\u00f75=\"],\n  [\"28\u00f79=\", \"51\u00f78=\"],\n  [\"77\u00f74=\", \"44\u00f74=\"],\n  [\"59\u00f74=\", \"67\u00f73=\"],\n  [\"78\u00f74=\", \"21\u00f72=\"],\n  [\"89\u00f75=\", \"47\u00f75=\"],\n  [\"98\u00f72=\", \"74\u00f79=\"],\n  [\"37\u00f78=\", \"44\u00f76=\"],\n  [\"35\u00f72=\", \"16\u00f74=\"],\n  [\"96\u00f78=\", \"75\u00f76=\"],\n  [\"32\u00f74=\", \"57\u00f73=\"],\n  [\"30\u00f75=\", \"39\u00f76=\"],\n  [\"43\u00f72=\", \"88\u00f79=\"],\n  [\"14\u00f78=\", \"59\u00f77=\"],\n  [\"82\u00f74=\", \"53\u00f74=\"],\n  [\"26\u00f79=\", \"19\u00f75=\"],\n  [\"24\u00f77=\", \"43\u00f75=\"],\n  [\"93\u00f74=\", \"49\u00f78=\"],\n  [\"69\u00f75=\", \"14\u00f77=\"],\n  [\"26\u00f73=\", \"14\u00f77=\"],\n  [\"97\u00f79=\", \"64\u00f72=\"],\n  [\"22\u00f77=\", \"75\u00f76=\"],\n  [\"43\u00f72=\", \"27\u00f79=\"],\n  [\"27\u00f75=\", \"65\u00f75=\"],\n];\n\nconst newValues = table.values.map((row) => row.slice());\nlet i = 0;\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const current = newValues[r][c];\n    if (current !== \"\") {\n      const [oldText, newText] = pairs[i];\n      if (current !== oldText) {\n        throw new Error(\n          `Unexpected cell text at row ${r}, col ${c}: expected \"${oldText}\", found \"${current}\"`\n        );\n      }\n      newValues[r][c] = newText;\n      i++;\n    }\n  }\n}\n\nif (i !== pairs.length) {\n  throw new Error(`Expected to replace ${pairs.length} cells, but replaced ${i}`);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the two-digit/one-digit division exercises in the first table.\n# Each data row of the table (rows 1, 5, 9, 13, 17 \u2014 the other rows are\n# blank \"answer\" rows) holds five \"NN\u00f7N=\" expressions; replace them in\n# reading order with their new values, matched positionally (not by text,\n# since some old values repeat with different replacements, e.g. \"43\u00f72=\"\n# appears twice and must become two different results).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# old/new pairs, in document (reading) order\n$oldTexts = @(\n    \"30\u00f74=\",\n    \"11\u00f77=\",\n    \"28\u00f79=\",\n    \"77\u00f74=\",\n    \"59\u00f74=\",\n    \"78\u00f74=\",\n    \"89\u00f75=\",\n    \"98\u00f72=\",\n    \"37\u00f78=\",\n    \"35\u00f72=\",\n    \"96\u00f78=\",\n    \"32\u00f74=\",\n    \"30\u00f75=\",\n    \"43\u00f72=\",\n    \"14\u00f78=\",\n    \"82\u00f74=\",\n    \"26\u00f79=\",\n    \"24\u00f77=\",\n    \"93\u00f74=\",\n    \"69\u00f75=\",\n    \"26\u00f73=\",\n    \"97\u00f79=\",\n    \"22\u00f77=\",\n    \"43\u00f72=\",\n    \"27\u00f75=\"\n)\n$newTexts = @(\n    \"45\u00f75=\",\n    \"56\u00f75=\",\n    \"51\u00f78=\",\n    \"44\u00f74=\",\n    \"67\u00f73=\",\n    \"21\u00f72=\",\n    \"47\u00f75=\",\n    \"74\u00f79=\",\n    \"44\u00f76=\",\n    \"16\u00f74=\",\n    \"75\u00f76=\",\n    \"57\u00f73=\",\n    \"39\u00f76=\",\n    \"88\u00f79=\",\n    \"59\u00f77=\",\n    \"53\u00f74=\",\n    \"19\u00f75=\",\n    \"43\u00f75=\",\n    \"49\u00f78=\",\n    \"14\u00f77=\",\n    \"14\u00f77=\",\n    \"64\u00f72=\",\n    \"75\u00f76=\",\n    \"27\u00f79=\",\n    \"65\u00f75=\"\n)\n\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $text = $cell.Range.Text\n        # Strip the trailing cell-mark characters (CR + BEL) that\n        # Range.Text includes for a table cell.\n        $plain = $text.TrimEnd([char]13, [char]7)\n        if ($plain -ne \"\") {\n            if ($plain -ne $oldTexts[$i]) {\n                throw \"Unexpected cell text at row $r, col $c`: expected '$($oldTexts[$i])', found '$plain'\"\n            }\n            $cell.Range.Text = $newTexts[$i]\n            $i = $i + 1\n        }\n    }\n}\n\nif ($i -ne $oldTexts.Count) {\n    throw \"Expected to replace $($oldTexts.Count) cells, but replaced $i\"\n}\n"}
